$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ row = 2; D = "U"; E = 1; F = 1; G = 1 },
    @{ row = 3; D = "U"; E = 1; F = 2; G = 2 },
    @{ row = 4; D = "U"; E = 1; F = 1; G = 2 },
    @{ row = 5; D = "U"; E = 1; F = 3; G = 3 }
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$ws.Range("D2:G5").Select()
